$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-6: increment the "Förändrad" (C) date value by 1 day
$ws.Range("C2").Value = 46060
$ws.Range("C3").Value = 46060
$ws.Range("C4").Value = 46060
$ws.Range("C5").Value = 46060
$ws.Range("C6").Value = 46060

# Rows 7-16: re-synced data (A, B, C, G columns) from source; row order also changed
$ws.Range("A7").Value = "A 32633-2025"
$ws.Range("B7").Value = 45838.65677083333
$ws.Range("C7").Value = 46060
$ws.Range("G7").Value = 1.3

$ws.Range("A8").Value = "A 45370-2022"
$ws.Range("B8").Value = 44844.6397337963
$ws.Range("C8").Value = 46060
$ws.Range("G8").Value = 2.7

$ws.Range("A9").Value = "A 23677-2023"
$ws.Range("B9").Value = 45077
$ws.Range("C9").Value = 46060
$ws.Range("G9").Value = 0.6

$ws.Range("A10").Value = "A 2253-2022"
$ws.Range("B10").Value = 44578
$ws.Range("C10").Value = 46060
$ws.Range("G10").Value = 0.3

$ws.Range("A11").Value = "A 58926-2025"
$ws.Range("B11").Value = 45986
$ws.Range("C11").Value = 46060
$ws.Range("G11").Value = 3.1

$ws.Range("A12").Value = "A 23678-2023"
$ws.Range("B12").Value = 45077
$ws.Range("C12").Value = 46060
$ws.Range("G12").Value = 1.4

$ws.Range("A13").Value = "A 50277-2024"
$ws.Range("B13").Value = 45600.60440972223
$ws.Range("C13").Value = 46060
$ws.Range("G13").Value = 0.5

$ws.Range("A14").Value = "A 50530-2024"
$ws.Range("B14").Value = 45601.56424768519
$ws.Range("C14").Value = 46060
$ws.Range("G14").Value = 0.7

$ws.Range("A15").Value = "A 50538-2024"
$ws.Range("B15").Value = 45601.57153935185
$ws.Range("C15").Value = 46060
$ws.Range("G15").Value = 0.8

$ws.Range("A16").Value = "A 11351-2021"
$ws.Range("B16").Value = 44263
$ws.Range("C16").Value = 46060
$ws.Range("G16").Value = 0.5

